$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Move "maquetear el proyecto parte visual" (B39) into C39, renamed to
# "maquetear el proyecto parte presupuesto" (task item moved to a different
# Kanban column with an updated title as part of the approved design).
$ws.Range("B39").Value = ""
$ws.Range("C39").Value = "maquetear el proyecto parte presupuesto"

# C40 keeps showing "Ajustar excel conciliacion para hacer calculos".
$ws.Range("C40").Value = "Ajustar excel conciliacion para hacer calculos"

# New pending task added in B41, with the row made taller to fit the
# wrapped two-line text.
$ws.Range("B41").Value = "Modificar generacion de presuúesto sin desplazamiento hacia abajo"
$ws.Rows.Item(41).RowHeight = 28.5

# The active selection now sits on the freshly added B41 cell.
$ws.Range("B41").Select()
